# Updated datum_excel_tests workbook in prep for unit test framework.
#
# - Rename Sheet1 -> "Gearbox Tests"
# - Add new named ranges for gearbox parameters (values live in column C)
# - Re-layout sheet1 with a small "gearbox parameters" table (labels in
#   column B, units in column D) and clear out the old scratch values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet first; existing defined names automatically track
# the new sheet name.
$ws.Name = "Gearbox Tests"

# --- Clear the old scratch contents -----------------------------------
# A1 gets overwritten below with a new string value.
# A2 ("float_range") becomes entirely empty again.
# A3 ("Date_range") keeps its date style/number format but loses its value.
$ws.Range("A2").Clear()
$ws.Range("A3").ClearContents()

# --- New cell content ---------------------------------------------------
# Order matters here: Excel assigns shared-string table indices in the
# order values are first written, so we write cells in the same order
# as the target file's shared string table.
$ws.Range("A1").Value = "Because spreadsheets need dark mode too!"

$ws.Range("B3").Value = "Gearbox parameters"

$ws.Range("B4").Value = "SURFACE_PAINTED.area"
$ws.Range("B5").Value = "HOUSING.mass"
$ws.Range("B6").Value = "FASTENERS.mass"
$ws.Range("B7").Value = "GEARS.mass"
$ws.Range("B8").Value = "DIPSTICK"
$ws.Range("B9").Value = "AIR NUT"
$ws.Range("B10").Value = "SHAFT CENTERS"

$ws.Range("D4").Value = "mm2"
$ws.Range("D5").Value = "kg"
$ws.Range("D6").Value = "kg"
$ws.Range("D7").Value = "kg"
$ws.Range("D8").Value = "deg"
$ws.Range("D10").Value = "mm"

# --- Column widths --------------------------------------------------------
$ws.Columns("B").ColumnWidth = 23.85546875

# --- Named ranges for the new gearbox parameter cells (column C) ----------
$wb.Names.Add("SURFACE_PAINTED.area", "='Gearbox Tests'!`$C`$4")
$wb.Names.Add("HOUSING.mass", "='Gearbox Tests'!`$C`$5")
$wb.Names.Add("FASTENERS.mass", "='Gearbox Tests'!`$C`$6")
$wb.Names.Add("GEARS.mass", "='Gearbox Tests'!`$C`$7")
$wb.Names.Add("DIPSTICK", "='Gearbox Tests'!`$C`$8")
$wb.Names.Add("AIR_NUT", "='Gearbox Tests'!`$C`$9")
$wb.Names.Add("SHAFT_CENTERS", "='Gearbox Tests'!`$C`$10")

# Leave the final selection on D10, matching the saved view state.
$ws.Range("D10").Select() | Out-Null
